$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/9/2024  Through  9/15/2024"

# --- Weekly crime-stat table updates (rows 15-30) ---
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = -47.826086956521
$ws.Range("N15").Value = -72.727272727272

# Row 16
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -46.153846153846
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = -21.551724137931
$ws.Range("L16").Value = 16.666666666666
$ws.Range("M16").Value = -55.392156862745
$ws.Range("N16").Value = -87.037037037037

# Row 17
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 3.125
$ws.Range("I17").Value = 345
$ws.Range("J17").Value = 328
$ws.Range("K17").Value = 5.182926829268
$ws.Range("L17").Value = 21.478873239436
$ws.Range("M17").Value = 61.971830985915
$ws.Range("N17").Value = -39.473684210526

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 1
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 900
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = -14.736842105263
$ws.Range("L18").Value = 10.958904109589
$ws.Range("M18").Value = -62.5
$ws.Range("N18").Value = -93.627065302911

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -25.714285714285
$ws.Range("I19").Value = 269
$ws.Range("J19").Value = 296
$ws.Range("K19").Value = -9.121621621621
$ws.Range("L19").Value = 8.467741935483
$ws.Range("M19").Value = -7.241379310344
$ws.Range("N19").Value = -33.250620347394

# Row 20
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 96
$ws.Range("J20").Value = 111
$ws.Range("K20").Value = -13.513513513513
$ws.Range("L20").Value = 20
$ws.Range("M20").Value = -32.394366197183
$ws.Range("N20").Value = -90

# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = -10.416666666666
$ws.Range("I21").Value = 895
$ws.Range("J21").Value = 968
$ws.Range("K21").Value = -7.541322314049
$ws.Range("L21").Value = 15.186615186615
$ws.Range("M21").Value = -18.413855970829
$ws.Range("N21").Value = -77.421796165489

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 800
$ws.Range("I23").Value = 61
$ws.Range("K23").Value = -1.612903225806
$ws.Range("L23").Value = 45.238095238095
$ws.Range("M23").Value = 74.285714285714

# Row 24
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -10.78431372549
$ws.Range("I24").Value = 887
$ws.Range("J24").Value = 886
$ws.Range("K24").Value = 0.112866817155
$ws.Range("L24").Value = 5.721096543504
$ws.Range("M24").Value = -17.334575955265

# Row 25
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 13.157894736842
$ws.Range("I25").Value = 394
$ws.Range("J25").Value = 342
$ws.Range("K25").Value = 15.204678362573
$ws.Range("L25").Value = 42.753623188405

# Row 26
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 66
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = 29.411764705882
$ws.Range("I26").Value = 556
$ws.Range("J26").Value = 512
$ws.Range("K26").Value = 8.59375
$ws.Range("L26").Value = 20.607375271149
$ws.Range("M26").Value = -37.52808988764

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("L27").Value = 4.166666666666

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -53.846153846153
$ws.Range("I28").Value = 65
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = -7.142857142857
$ws.Range("L28").Value = 20.37037037037

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("L29").Value = -65

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "***.*"
$ws.Range("L30").Value = -53.333333333333

